# Adds 4 new flight-departure rows (159-162) to the "Main Data" sheet,
# continuing the existing table that ends at row 158.
# simple graph added, requires small fixes. added todo tasks

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 159; A = 158; B = "Saturday, Jan 14"; C = "6:35 PM";  D = "FR7946"; E = "Oslo";      F = "(TRF)"; G = "Ryanair "; H = "B738"; I = "(SP-RSX)"; J = "6:37 PM";  L = "0 hours, 2 minutes" },
    @{ Row = 160; A = 159; B = "Saturday, Jan 14"; C = "7:10 PM";  D = "FR4272"; E = "Budapest";  F = "(BUD)"; G = "Ryanair "; H = "B738"; I = "(SP-RSM)"; J = "7:35 PM";  L = "0 hours, 25 minutes" },
    @{ Row = 161; A = 160; B = "Saturday, Jan 14"; C = "9:55 PM";  D = "FR1751"; E = "London";    F = "(STN)"; G = "Ryanair "; H = "B38M"; I = "(EI-HEV)"; J = "10:15 PM"; L = "0 hours, 20 minutes" },
    @{ Row = 162; A = 161; B = "Saturday, Jan 14"; C = "10:30 PM"; D = "FR9629"; E = "Liverpool"; F = "(LPL)"; G = "Ryanair "; H = "B738"; I = "(EI-EKZ)"; J = "10:49 PM"; L = "0 hours, 19 minutes" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $r.A
    $ws.Cells.Item($row, 2).Value  = $r.B
    $ws.Cells.Item($row, 3).Value  = $r.C
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 5).Value  = $r.E
    $ws.Cells.Item($row, 6).Value  = $r.F
    $ws.Cells.Item($row, 7).Value  = $r.G
    $ws.Cells.Item($row, 8).Value  = $r.H
    $ws.Cells.Item($row, 9).Value  = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    # Column K (DIFFERENCE) is intentionally left blank, matching the rest of the table.
    $ws.Cells.Item($row, 12).Value = $r.L
    # Column M is intentionally left blank, matching the rest of the table.
}
